$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered) from O1 to the new P1:Q1 header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Set new header values in P1 and Q1
$ws.Cells.Item(1, 16).Value2 = 14
$ws.Cells.Item(1, 17).Value2 = 15

# For rows 2-25: swap I/K values, swap M/O values, and add P/Q = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 9).Value2 = $kVal
    $ws.Cells.Item($r, 11).Value2 = $iVal

    $mVal = $ws.Cells.Item($r, 13).Value2
    $oVal = $ws.Cells.Item($r, 15).Value2
    $ws.Cells.Item($r, 13).Value2 = $oVal
    $ws.Cells.Item($r, 15).Value2 = $mVal

    $ws.Cells.Item($r, 16).Value2 = 2
    $ws.Cells.Item($r, 17).Value2 = 2
}
